# Weekly update: insert a new daily price record for
# "Feria Lagunitas de Puerto Montt" / Mandarina / Clementina / Primera,
# shifting the existing rows 63-98 down to 64-99 (their data is untouched)
# and filling the newly opened row 63 with the new week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 63; Excel shifts rows 63:98 -> 64:99
# and the sheet's used range grows from A1:T98 to A1:T99 automatically.
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with the new observation.
$ws.Cells.Item(63, 1).Value  = 4
$ws.Cells.Item(63, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(63, 3).Value  = "Los Lagos"
$ws.Cells.Item(63, 4).Value  = 44452
$ws.Cells.Item(63, 5).Value  = 10
$ws.Cells.Item(63, 6).Value  = "Fruta"
$ws.Cells.Item(63, 7).Value  = 100102
$ws.Cells.Item(63, 8).Value  = "Cítricos"
$ws.Cells.Item(63, 9).Value  = 100102004
$ws.Cells.Item(63, 10).Value = "Mandarina"
$ws.Cells.Item(63, 11).Value = "Clementina"
$ws.Cells.Item(63, 12).Value = "Primera"
$ws.Cells.Item(63, 13).Value = 300
$ws.Cells.Item(63, 14).Value = 6500
$ws.Cells.Item(63, 15).Value = 6500
$ws.Cells.Item(63, 16).Value = 6500
$ws.Cells.Item(63, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(63, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(63, 19).Value = 650
$ws.Cells.Item(63, 20).Value = 10
